$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
# Row 17
$ws.Range("H17").Value = 4008327.5
$ws.Range("J17").Value = 4175262
$ws.Range("L17").Value = 12525786
$ws.Range("N17").Value = -12526122
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
# Row 38
$ws.Range("H38").Value = 329.7143
$ws.Range("I38").Value = 57.333332
$ws.Range("J38").Value = 820
$ws.Range("K38").Value = 171.999996
$ws.Range("L38").Value = 2460
$ws.Range("M38").Value = 200.000004
$ws.Range("N38").Value = -3204
# Row 51
$ws.Range("H51").Value = 2244.4443
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2244.4443
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2244.4443
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -3212.4443
# Row 58
$ws.Range("H58").Value = 748.3333
$ws.Range("I58").Value = 676.4286
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 2029.2858
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -1879.2858
$ws.Range("N58").Value = -3300
# Row 76
$ws.Range("H76").Value = 1918852.6
$ws.Range("I76").Value = 3212.2273
$ws.Range("J76").Value = 7939436.5
$ws.Range("K76").Value = 3212.2273
$ws.Range("L76").Value = 7939436.5
$ws.Range("M76").Value = -2897.2273
$ws.Range("N76").Value = -7940066.5
# Row 79
$ws.Range("H79").Value = 1918852.6
$ws.Range("I79").Value = 3212.2273
$ws.Range("J79").Value = 7939436.5
$ws.Range("K79").Value = 3212.2273
$ws.Range("L79").Value = 7939436.5
$ws.Range("M79").Value = -2120.2273
$ws.Range("N79").Value = -7941620.5
# Row 98
$ws.Range("H98").Value = 661.875
$ws.Range("I98").Value = 670.7143
$ws.Range("J98").Value = 600
$ws.Range("K98").Value = 670.7143
$ws.Range("L98").Value = 600
$ws.Range("M98").Value = 827.2857
$ws.Range("N98").Value = -3596
# Row 122
$ws.Range("H122").Value = 661.875
$ws.Range("I122").Value = 670.7143
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 2012.1429
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = 437.8571000000002
$ws.Range("N122").Value = -6700

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1568.1578
$ws.Range("I2").Value = 1406.2941
$ws.Range("K2").Value = 1406.2941
$ws.Range("M2").Value = -1293.2941
# Row 23
$ws.Range("H23").Value = 3835.3333
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 3835.3333
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 3835.3333
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -4353.3333
# Row 32
$ws.Range("H32").Value = 5353.6665
$ws.Range("I32").Value = 5041.137
$ws.Range("K32").Value = 5041.137
$ws.Range("M32").Value = -4754.137
# Row 45
$ws.Range("H45").Value = 2987.7693
$ws.Range("I45").Value = 3199.8
$ws.Range("J45").Value = 2855.25
$ws.Range("K45").Value = 3199.8
$ws.Range("L45").Value = 2855.25
$ws.Range("M45").Value = -2822.8
$ws.Range("N45").Value = -3609.25
# Row 63
$ws.Range("H63").Value = 2842256
$ws.Range("I63").Value = 1481.5
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 1481.5
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -795.5
$ws.Range("N63").Value = -31251372
# Row 66
$ws.Range("H66").Value = 2842256
$ws.Range("I66").Value = 1481.5
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 7407.5
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -3975.5
$ws.Range("N66").Value = -156256864
# Row 116
$ws.Range("H116").Value = 1568.1578
$ws.Range("I116").Value = 1406.2941
$ws.Range("K116").Value = 1406.2941
$ws.Range("M116").Value = 887.7058999999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1568.1578
$ws.Range("I3").Value = 1406.2941
$ws.Range("K3").Value = 1406.2941
$ws.Range("M3").Value = -1292.2941
# Row 29
$ws.Range("H29").Value = 500
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
# Row 35
$ws.Range("H35").Value = 22794.8
$ws.Range("J35").Value = 22794.8
$ws.Range("L35").Value = 22794.8
$ws.Range("N35").Value = -23414.8
# Row 82
$ws.Range("H82").Value = 23562
$ws.Range("I82").Value = 7974.75
$ws.Range("J82").Value = 44345
$ws.Range("K82").Value = 7974.75
$ws.Range("L82").Value = 44345
$ws.Range("M82").Value = -7591.75
$ws.Range("N82").Value = -45111
# Row 85
$ws.Range("H85").Value = 23562
$ws.Range("I85").Value = 7974.75
$ws.Range("J85").Value = 44345
$ws.Range("K85").Value = 7974.75
$ws.Range("L85").Value = 44345
$ws.Range("M85").Value = -6648.75
$ws.Range("N85").Value = -46997
# Row 94
$ws.Range("H94").Value = 1400
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1400
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1400
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2302

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 47
$ws.Range("H47").Value = 13833.333
$ws.Range("I47").Value = 7500
$ws.Range("J47").Value = 17000
$ws.Range("K47").Value = 7500
$ws.Range("L47").Value = 17000
$ws.Range("M47").Value = -6934
$ws.Range("N47").Value = -18132
# Row 62
$ws.Range("H62").Value = 4353.25
$ws.Range("I62").Value = 4676.143
$ws.Range("J62").Value = 3901.2
$ws.Range("K62").Value = 4676.143
$ws.Range("L62").Value = 3901.2
$ws.Range("M62").Value = -4052.143
$ws.Range("N62").Value = -5149.2
# Row 65
$ws.Range("H65").Value = 4353.25
$ws.Range("I65").Value = 4676.143
$ws.Range("J65").Value = 3901.2
$ws.Range("K65").Value = 23380.715
$ws.Range("L65").Value = 19506
$ws.Range("M65").Value = -20260.715
$ws.Range("N65").Value = -25746
# Row 94
$ws.Range("H94").Value = 1410.5454
$ws.Range("J94").Value = 2486
$ws.Range("L94").Value = 2486
$ws.Range("N94").Value = -3388
# Row 132
$ws.Range("H132").Value = 3611.2
$ws.Range("I132").Value = 2460.8
$ws.Range("K132").Value = 7382.400000000001
$ws.Range("M132").Value = -4852.400000000001
# Row 134
$ws.Range("H134").Value = 1284.2609
$ws.Range("I134").Value = 1209.8572
$ws.Range("J134").Value = 1400
$ws.Range("K134").Value = 3629.5716
$ws.Range("L134").Value = 4200
$ws.Range("M134").Value = -1094.5716
$ws.Range("N134").Value = -9270

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1557.4
$ws.Range("I5").Value = 1551.1428
$ws.Range("J5").Value = 1560.7693
$ws.Range("K5").Value = 4653.428400000001
$ws.Range("L5").Value = 4682.3079
$ws.Range("M5").Value = -4541.428400000001
$ws.Range("N5").Value = -4906.3079
# Row 64
$ws.Range("H64").Value = 2477.1428
$ws.Range("I64").Value = 1406
$ws.Range("J64").Value = 2905.6
$ws.Range("K64").Value = 4218
$ws.Range("L64").Value = 8716.799999999999
$ws.Range("M64").Value = -3948
$ws.Range("N64").Value = -9256.799999999999
# Row 67
$ws.Range("H67").Value = 2477.1428
$ws.Range("I67").Value = 1406
$ws.Range("J67").Value = 2905.6
$ws.Range("K67").Value = 4218
$ws.Range("L67").Value = 8716.799999999999
$ws.Range("M67").Value = -3282
$ws.Range("N67").Value = -10588.8
# Row 75
$ws.Range("H75").Value = 761.6667
$ws.Range("I75").Value = 761.6667
$ws.Range("K75").Value = 2285.0001
$ws.Range("M75").Value = -1287.0001
# Row 78
$ws.Range("H78").Value = 761.6667
$ws.Range("I78").Value = 761.6667
$ws.Range("K78").Value = 6855.0003
$ws.Range("M78").Value = -1863.0003
# Row 131
$ws.Range("H131").Value = 730.77
$ws.Range("I131").Value = 546
$ws.Range("J131").Value = 740.49475
$ws.Range("K131").Value = 1638
$ws.Range("L131").Value = 2221.48425
$ws.Range("M131").Value = 3402
$ws.Range("N131").Value = -12301.48425
# Row 135
$ws.Range("H135").Value = 1557.4
$ws.Range("I135").Value = 1551.1428
$ws.Range("J135").Value = 1560.7693
$ws.Range("K135").Value = 13960.2852
$ws.Range("L135").Value = 14046.9237
$ws.Range("M135").Value = -11425.2852
$ws.Range("N135").Value = -19116.9237

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 5350.1934
$ws.Range("J126").Value = 6453.143
$ws.Range("L126").Value = 19359.429
$ws.Range("N126").Value = -24299.429

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 35
$ws.Range("H35").Value = 251157.25
$ws.Range("I35").Value = 251157.25
$ws.Range("K35").Value = 251157.25
$ws.Range("M35").Value = -250821.25
# Row 40
$ws.Range("H40").Value = 3075.7878
$ws.Range("I40").Value = 2452.074
$ws.Range("K40").Value = 2452.074
$ws.Range("M40").Value = -2316.074
# Row 122
$ws.Range("H122").Value = 702736.1
$ws.Range("I122").Value = 1636029
$ws.Range("J122").Value = 2766.5
$ws.Range("K122").Value = 4908087
$ws.Range("L122").Value = 8299.5
$ws.Range("M122").Value = -4905637
$ws.Range("N122").Value = -13199.5
# Row 132
$ws.Range("H132").Value = 3680.9
$ws.Range("J132").Value = 5332.3335
$ws.Range("L132").Value = 15997.0005
$ws.Range("N132").Value = -21057.0005
# Row 136
$ws.Range("H136").Value = 1835.8667
$ws.Range("J136").Value = 1350
$ws.Range("L136").Value = 4050
$ws.Range("N136").Value = -9150

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 15153281
$ws.Range("I107").Value = 2799
$ws.Range("J107").Value = 22728522
$ws.Range("K107").Value = 8397
$ws.Range("L107").Value = 68185566
$ws.Range("M107").Value = -6477
$ws.Range("N107").Value = -68189406
# Row 122
$ws.Range("H122").Value = 1002.5238
$ws.Range("I122").Value = 959.5
$ws.Range("J122").Value = 1088.5714
$ws.Range("K122").Value = 2878.5
$ws.Range("L122").Value = 3265.7142
$ws.Range("M122").Value = -428.5
$ws.Range("N122").Value = -8165.7142
